# Add 2022-Q3 sheet (new fund-holding snapshot) and keep the existing
# 2021-Q2 sheet, updating the "总计" (totals) summary sheet accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by duplicating the current "2021-Q2"
#    sheet (so it inherits the same column formatting/borders), then
#    rename sheets so sheetId allocation matches: 2022-Q3 keeps the
#    original sheetId, 2021-Q2 (the untouched duplicate) gets the new one.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($null, $q2Sheet)

$newQ3 = $wb.Worksheets.Item(2)
$keepQ2 = $wb.Worksheets.Item(3)

$newQ3.Name = "2022-Q3"
$keepQ2.Name = "2021-Q2"

# ---------------------------------------------------------------------
# 2) Populate "2022-Q3" with the new fund-holding data.
# ---------------------------------------------------------------------

# Header D1 differs slightly in wording for the new sheet.
$newQ3.Range("D1").Value = "基金规模"

# Extend the index-column (A) formatting down to the two extra rows.
$newQ3.Range("A2").Copy()
$newQ3.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the old (2021-Q2) data from rows 2-3; rows 4-5 are still empty.
$newQ3.Rows("2:3").ClearContents()

# Fund code / name / size / position columns must stay text (e.g. keep
# leading zeros in fund codes, and "1.77" rather than 1.77).
$newQ3.Range("B2:G5").NumberFormat = "@"

$q3Data = @(
    @("014575", "鑫元清洁能源混合C", "1.77", "93.05", "7.25", "0.1283", 8),
    @("014574", "鑫元清洁能源混合A", "0.82", "93.05", "7.25", "0.0594", 8),
    @("013470", "泰信低碳经济混合C", "0.24", "82.26", "4.11", "0.0099", 6),
    @("013469", "泰信低碳经济混合A", "0.16", "82.26", "4.11", "0.0066", 6)
)

$r = 2
foreach ($row in $q3Data) {
    $newQ3.Cells.Item($r, 1).Value = $r - 2
    $newQ3.Cells.Item($r, 2).Value = $row[0]
    $newQ3.Cells.Item($r, 3).Value = $row[1]
    $newQ3.Cells.Item($r, 4).Value = $row[2]
    $newQ3.Cells.Item($r, 5).Value = $row[3]
    $newQ3.Cells.Item($r, 6).Value = $row[4]
    $newQ3.Cells.Item($r, 7).Value = $row[5]
    $newQ3.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# Match the page margins used elsewhere in the workbook.
$newQ3.PageSetup.LeftMargin = $excel.Application.InchesToPoints(0.75)
$newQ3.PageSetup.RightMargin = $excel.Application.InchesToPoints(0.75)
$newQ3.PageSetup.TopMargin = $excel.Application.InchesToPoints(1)
$newQ3.PageSetup.BottomMargin = $excel.Application.InchesToPoints(1)
$newQ3.PageSetup.HeaderMargin = $excel.Application.InchesToPoints(0.5)
$newQ3.PageSetup.FooterMargin = $excel.Application.InchesToPoints(0.5)

# ---------------------------------------------------------------------
# 3) Update the "总计" (totals) sheet: insert the 2022-Q3 totals above
#    the existing 2021-Q2 row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# Give the (new) row 3 index cell the same style as row 2's index cell.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the existing 2021-Q2 totals down to row 3 (index becomes 1).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.01

# Write the new 2022-Q3 totals into row 2 (A2 already holds index 0).
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.2
